$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Add a new user-story bullet paragraph after "As a user, I want to
#    be able to schedule for more employees in response to the
#    weather." — matching the same ListParagraph / numId=2 bullet
#    style, with the text split across three runs.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*schedule for more employees in response to the weather*") {
        $target = $cand
        break
    }
}

$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($target.Index + 1)

$newParaXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">As a user, I want to be able </w:t></w:r><w:r><w:t xml:space="preserve">to </w:t></w:r><w:r><w:t>easily manipulate the program by myself.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newParaRange.InsertXML($newParaXml)

# ---------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker: it currently sits on
#    the "Number of seniors" run, but it should sit on the "Number of
#    customers for dinner rush" run (the previous list item).
# ---------------------------------------------------------------------
$dinnerPara = $null
$seniorsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t -like "*Number of customers for dinner rush*") {
        $dinnerPara = $cand
    }
    if ($t -like "*Number of seniors*") {
        $seniorsPara = $cand
    }
}

$dinnerRange = $d.Range($dinnerPara.Range.Start, $dinnerPara.Range.End - 1)
$dinnerXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Number of customers for dinner rush</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$dinnerRange.InsertXML($dinnerXml)

$seniorsRange = $d.Range($seniorsPara.Range.Start, $seniorsPara.Range.End - 1)
$seniorsXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Number of seniors</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$seniorsRange.InsertXML($seniorsXml)

# ---------------------------------------------------------------------
# 3) Remove the trailing empty ListParagraph (ilvl=1, numId=1) that
#    follows the "Weather is currently defined..." paragraph, right
#    before the final sectPr.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$veryLast = $d.Paragraphs.Item($count)
if ($veryLast.Range.Text -eq "") {
    $prev = $d.Paragraphs.Item($count - 1)
    $delRange = $d.Range($prev.Range.End - 1, $veryLast.Range.End)
    $delRange.Delete()
}

Write-Output "done"
